$d = $word.ActiveDocument

# --- Paragraph 1: "This checkbox will be checked: " + checked-box glyph (MS Gothic) ---
$para1Xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">This checkbox will be checked: </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" w:ascii="MS Gothic" w:hAnsi="MS Gothic" w:eastAsia="MS Gothic"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t>&#x2612;</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

# --- Paragraph 2: "This checkbox will be " + "un" + "checked: " + unchecked-box glyph (MS Gothic) ---
$para2Xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">This checkbox will be </w:t></w:r><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t>un</w:t></w:r><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">checked: </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" w:ascii="MS Gothic" w:hAnsi="MS Gothic" w:eastAsia="MS Gothic"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t>&#x2610;</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

# The document starts as a single empty paragraph. InsertXML on that (still
# empty) paragraph's Range inserts a brand-new paragraph immediately before
# it, leaving the original (still-empty) paragraph intact right after it.
# Doing this twice - once per target paragraph - builds both checkbox lines
# while leaving one now-superfluous empty paragraph behind at the end.
$originalPara = $d.Paragraphs(1)
$originalPara.Range.InsertXML($para1Xml)

$originalPara = $d.Paragraphs(2)
$originalPara.Range.InsertXML($para2Xml)

# Clean up the leftover empty trailing paragraph (merge its paragraph mark
# into paragraph 2 by deleting from the end of paragraph 2's text through
# the end of the document).
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$trailingRange = $d.Range($d.Paragraphs(2).Range.End - 1, $lastPara.Range.End)
$trailingRange.Delete()

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
